# Generate Report for Handoff
# The file f04ac914-b1e8-41a4-91d7-d843ccb1764a has finished translation and is
# now ready for handoff. Update its status/priority/handoff-timestamp on every
# sheet of the localization status report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (row 3 = f04ac914-...) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-30 14:16:04"
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet (row 3 = f04ac914-...) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-30 14:15:56"
$zhcn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet (row 3 = f04ac914-...) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-30 14:16:04"
$dede.Columns.Item(3).ColumnWidth = 16.3
